$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 412; this shifts the existing rows 412:473
# down to 416:477 and updates the sheet dimension automatically.
$ws.Rows("412:415").Insert()

# New data block (newest price report) that now occupies rows 412:415.
$data = @(
    @{ Row = 412; Calidad = "Especial"; Volumen = 300; Min = 11000; Max = 12000; Prom = 11500; PrecioKg = 639 },
    @{ Row = 413; Calidad = "Primera";  Volumen = 500; Min = 9000;  Max = 10000; Prom = 9500;  PrecioKg = 528 },
    @{ Row = 414; Calidad = "Segunda";  Volumen = 400; Min = 7000;  Max = 8000;  Prom = 7500;  PrecioKg = 417 },
    @{ Row = 415; Calidad = "Tercera";  Volumen = 340; Min = 4000;  Max = 5000;  Prom = 4500;  PrecioKg = 250 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = 2
    $ws.Cells.Item($r, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = 44776
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = 100112043
    $ws.Cells.Item($r, 7).Value = "Pepino dulce"
    $ws.Cells.Item($r, 8).Value = "Cultivar IV Región"
    $ws.Cells.Item($r, 9).Value = $item.Calidad
    $ws.Cells.Item($r, 10).Value = $item.Volumen
    $ws.Cells.Item($r, 11).Value = $item.Min
    $ws.Cells.Item($r, 12).Value = $item.Max
    $ws.Cells.Item($r, 13).Value = $item.Prom
    $ws.Cells.Item($r, 14).Value = "$/bandeja 18 kilos"
    $ws.Cells.Item($r, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item($r, 16).Value = $item.PrecioKg
    $ws.Cells.Item($r, 17).Value = 18
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
